$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3818.375
$ws.Range("I137").Value = 1051.2727
$ws.Range("J137").Value = 4867.9653
$ws.Range("K137").Value = 3153.8181
$ws.Range("L137").Value = 14603.8959
$ws.Range("M137").Value = -603.8181
$ws.Range("N137").Value = -19703.8959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3072.7856
$ws.Range("I61").Value = 1601.6666
$ws.Range("J61").Value = 4176.125
$ws.Range("K61").Value = 1601.6666
$ws.Range("L61").Value = 4176.125
$ws.Range("M61").Value = -1389.6666
$ws.Range("N61").Value = -4600.125
$ws.Range("H74").Value = 2210.1724
$ws.Range("I74").Value = 1932.4762
$ws.Range("J74").Value = 2939.125
$ws.Range("K74").Value = 1932.4762
$ws.Range("L74").Value = 2939.125
$ws.Range("M74").Value = -1058.4762
$ws.Range("N74").Value = -4687.125
$ws.Range("H77").Value = 2210.1724
$ws.Range("I77").Value = 1932.4762
$ws.Range("J77").Value = 2939.125
$ws.Range("K77").Value = 9662.381000000001
$ws.Range("L77").Value = 14695.625
$ws.Range("M77").Value = -5294.381000000001
$ws.Range("N77").Value = -23431.625
$ws.Range("H122").Value = 2203.2
$ws.Range("I122").Value = 2026
$ws.Range("K122").Value = 6078
$ws.Range("M122").Value = -3628
$ws.Range("H136").Value = 3072.7856
$ws.Range("I136").Value = 1601.6666
$ws.Range("J136").Value = 4176.125
$ws.Range("K136").Value = 4804.9998
$ws.Range("L136").Value = 12528.375
$ws.Range("M136").Value = -2254.9998
$ws.Range("N136").Value = -17628.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 35000
$ws.Range("I15").Value = 20000
$ws.Range("J15").Value = 40000
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 40000
$ws.Range("M15").Value = -19773
$ws.Range("N15").Value = -40454
$ws.Range("H74").Value = 13244.75
$ws.Range("J74").Value = 13244.75
$ws.Range("L74").Value = 13244.75
$ws.Range("N74").Value = -15116.75
$ws.Range("H77").Value = 13244.75
$ws.Range("J77").Value = 13244.75
$ws.Range("L77").Value = 39734.25
$ws.Range("N77").Value = -49094.25
$ws.Range("H81").Value = 21479.8
$ws.Range("J81").Value = 21479.8
$ws.Range("L81").Value = 21479.8
$ws.Range("N81").Value = -23601.8
$ws.Range("H84").Value = 21479.8
$ws.Range("J84").Value = 21479.8
$ws.Range("L84").Value = 64439.39999999999
$ws.Range("N84").Value = -75047.39999999999
$ws.Range("H111").Value = 40002
$ws.Range("J111").Value = 40002
$ws.Range("L111").Value = 40002
$ws.Range("N111").Value = -48182
$ws.Range("H134").Value = 2481.7163
$ws.Range("I134").Value = 1462.0541
$ws.Range("K134").Value = 4386.1623
$ws.Range("M134").Value = -1851.1623

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6178460
$ws.Range("I31").Value = 3213.25
$ws.Range("J31").Value = 8778564
$ws.Range("K31").Value = 3213.25
$ws.Range("L31").Value = 8778564
$ws.Range("M31").Value = -2918.25
$ws.Range("N31").Value = -8779154
$ws.Range("H34").Value = 6178460
$ws.Range("I34").Value = 3213.25
$ws.Range("J34").Value = 8778564
$ws.Range("K34").Value = 3213.25
$ws.Range("L34").Value = 8778564
$ws.Range("M34").Value = -3011.25
$ws.Range("N34").Value = -8778968
$ws.Range("H58").Value = 2183.4827
$ws.Range("I58").Value = 1571.9375
$ws.Range("J58").Value = 2936.1538
$ws.Range("K58").Value = 1571.9375
$ws.Range("L58").Value = 2936.1538
$ws.Range("M58").Value = -1368.9375
$ws.Range("N58").Value = -3342.1538
$ws.Range("H132").Value = 44424.12
$ws.Range("I132").Value = 1423.0454
$ws.Range("J132").Value = 130426.27
$ws.Range("K132").Value = 4269.1362
$ws.Range("L132").Value = 391278.81
$ws.Range("M132").Value = -1739.1362
$ws.Range("N132").Value = -396338.81
$ws.Range("H134").Value = 45629.906
$ws.Range("I134").Value = 1483.35
$ws.Range("K134").Value = 4450.049999999999
$ws.Range("M134").Value = -1915.049999999999
$ws.Range("H136").Value = 2183.4827
$ws.Range("I136").Value = 1571.9375
$ws.Range("J136").Value = 2936.1538
$ws.Range("K136").Value = 4715.8125
$ws.Range("L136").Value = 8808.4614
$ws.Range("M136").Value = -2165.8125
$ws.Range("N136").Value = -13908.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -9707
$ws.Range("N18").Value = -10586
$ws.Range("H126").Value = 9697.357
$ws.Range("I126").Value = 19408.166
$ws.Range("J126").Value = 2414.25
$ws.Range("K126").Value = 58224.49800000001
$ws.Range("L126").Value = 7242.75
$ws.Range("M126").Value = -55754.49800000001
$ws.Range("N126").Value = -12182.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3495.9167
$ws.Range("I132").Value = 3202.56
$ws.Range("J132").Value = 4162.636
$ws.Range("K132").Value = 9607.68
$ws.Range("L132").Value = 12487.908
$ws.Range("M132").Value = -7077.68
$ws.Range("N132").Value = -17547.908
$ws.Range("H136").Value = 1775.7188
$ws.Range("I136").Value = 1408.88
$ws.Range("J136").Value = 3085.8572
$ws.Range("K136").Value = 4226.64
$ws.Range("L136").Value = 9257.571599999999
$ws.Range("M136").Value = -1676.64
$ws.Range("N136").Value = -14357.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 40406.5
$ws.Range("J108").Value = 40406.5
$ws.Range("L108").Value = 40406.5
$ws.Range("N108").Value = -48086.5
$ws.Range("H132").Value = 2268.8462
$ws.Range("I132").Value = 1606.5333
$ws.Range("K132").Value = 4819.5999
$ws.Range("M132").Value = -2289.5999
$ws.Range("H136").Value = 179467.58
$ws.Range("I136").Value = 257156.72
$ws.Range("J136").Value = 1239.5294
$ws.Range("K136").Value = 771470.16
$ws.Range("L136").Value = 3718.5882
$ws.Range("M136").Value = -768920.16
$ws.Range("N136").Value = -8818.5882
